$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Set A5 to the new value "Sensor Ultrasonico"
$ws.Range("A5").Value = "Sensor Ultrasonico"

# Widen column A (target raw width 17.453125; ColumnWidth is stored/quantized
# in 1/6-character steps by this runtime, so feed it the equivalent
# character-width value - i.e. raw width minus the 5/6 character padding)
$ws.Columns.Item(1).ColumnWidth = 17.453125 - (5/6)

# Move active selection to A6 (mirrors typing + Enter moving down)
$ws.Range("A6").Select()
